$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $value) {
    # Force the cell to remain plain text even when the new value looks
    # like a number (e.g. "59.22") or a signed percentage, matching the
    # inlineStr / shared-string cells openpyxl wrote originally.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextCell $ws.Range('D2') '43.851.88'
Set-TextCell $ws.Range('E2') '  -0.43%  '
Set-TextCell $ws.Range('D3') '2.341.48'
Set-TextCell $ws.Range('E3') '  -0.61%  '
Set-TextCell $ws.Range('E4') '  +0.04%  '
Set-TextCell $ws.Range('D5') '239.06'
Set-TextCell $ws.Range('E5') '  -1.16%  '
Set-TextCell $ws.Range('D6') '0.665'
Set-TextCell $ws.Range('E6') '  -4.30%  '
Set-TextCell $ws.Range('D7') '72.49'
Set-TextCell $ws.Range('E7') '  -5.33%  '
Set-TextCell $ws.Range('D9') '0.593'
Set-TextCell $ws.Range('E9') '  -6.53%  '
Set-TextCell $ws.Range('E10') '  -1.39%  '
Set-TextCell $ws.Range('D11') '59.22'
Set-TextCell $ws.Range('E11') '  +3.17%  '
Set-TextCell $ws.Range('D12') '32.68'
Set-TextCell $ws.Range('E12') '  -2.43%  '
Set-TextCell $ws.Range('E13') '  +0.18%  '
Set-TextCell $ws.Range('E14') '  -4.06%  '
Set-TextCell $ws.Range('D15') '2.692.68'
Set-TextCell $ws.Range('E15') '  -0.47%  '
Set-TextCell $ws.Range('D16') '16.03'
Set-TextCell $ws.Range('E16') '  -4.87%  '
Set-TextCell $ws.Range('D17') '0.896'
Set-TextCell $ws.Range('E17') '  -3.90%  '
Set-TextCell $ws.Range('D18') '2.337.50'
Set-TextCell $ws.Range('E18') '  -0.61%  '
Set-TextCell $ws.Range('D19') '43.758.45'
Set-TextCell $ws.Range('E19') '  -0.25%  '
Set-TextCell $ws.Range('E20') '  +0.14%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextCell $ws.Range('D21') '6.65'
Set-TextCell $ws.Range('E21') '  -0.36%  '
$ws.Range('B22').Value = 'Litecoin'
$ws.Range('C22').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell $ws.Range('D22') '78.19'
Set-TextCell $ws.Range('E22') '  +0.54%  '
Set-TextCell $ws.Range('D23') '250.93'
Set-TextCell $ws.Range('E23') '  -4.47%  '
Set-TextCell $ws.Range('E24') '  +0.16%  '
Set-TextCell $ws.Range('D25') '3.74'
Set-TextCell $ws.Range('E25') '  +3.02%  '
Set-TextCell $ws.Range('E26') '  +1.84%  '
Set-TextCell $ws.Range('D27') '2.48'
Set-TextCell $ws.Range('E27') '  -2.08%  '
Set-TextCell $ws.Range('D28') '10.38'
Set-TextCell $ws.Range('E28') '  -5.51%  '
Set-TextCell $ws.Range('E29') '  -1.81%  '
Set-TextCell $ws.Range('D30') '176.76'
Set-TextCell $ws.Range('E30') '  +0.81%  '
Set-TextCell $ws.Range('D31') '22.16'
Set-TextCell $ws.Range('E31') '  -4.26%  '
Set-TextCell $ws.Range('E32') '  -1.04%  '
Set-TextCell $ws.Range('E33') '  -3.14%  '
Set-TextCell $ws.Range('D34') '0.0745'
Set-TextCell $ws.Range('E34') '  -2.60%  '
Set-TextCell $ws.Range('D35') '5.06'
Set-TextCell $ws.Range('E35') '  -6.28%  '
Set-TextCell $ws.Range('D36') '5.33'
Set-TextCell $ws.Range('E36') '  -1.85%  '
Set-TextCell $ws.Range('E37') '  -2.67%  '
Set-TextCell $ws.Range('D38') '6.38'
Set-TextCell $ws.Range('E38') '  -0.64%  '
$ws.Range('B39').Value = 'FTXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextCell $ws.Range('D39') '5.81'
Set-TextCell $ws.Range('E39') '  +23.88%  '
$ws.Range('B40').Value = 'LidoDAOToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell $ws.Range('D40') '2.36'
Set-TextCell $ws.Range('E40') '  -2.48%  '
Set-TextCell $ws.Range('E41') '  -4.47%  '
Set-TextCell $ws.Range('D42') '65.35'
Set-TextCell $ws.Range('E42') '  +15.34%  '
Set-TextCell $ws.Range('D43') '9.21'
Set-TextCell $ws.Range('E43') '  +0.48%  '
Set-TextCell $ws.Range('E44') '  -1.16%  '
Set-TextCell $ws.Range('D45') '18.73'
Set-TextCell $ws.Range('E45') '  -2.88%  '
Set-TextCell $ws.Range('D46') '0.194'
Set-TextCell $ws.Range('E46') '  -12.39%  '
Set-TextCell $ws.Range('E47') '  +0.13%  '
Set-TextCell $ws.Range('E48') '  -3.14%  '
Set-TextCell $ws.Range('D49') '1.14'
Set-TextCell $ws.Range('E49') '  -3.90%  '
Set-TextCell $ws.Range('E50') '  -5.56%  '
Set-TextCell $ws.Range('D51') '97.83'
Set-TextCell $ws.Range('E51') '  -4.16%  '
